# Apply updated NATMI ligand-receptor pair statistics (F8 -> Lrp1)
# Following Dr Hou's advice: ligand/receptor expressing cell counts increase
# from 1 to 3, and all downstream derived expression/specificity metrics are
# recomputed accordingly. Values below are taken verbatim from the recomputed
# NATMI output table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2"=3; "G2"=2.229372666666666; "H2"=6.688117999999999; "I2"=0.5889623983027473; "J2"=0.5889623983027473; "K2"=3; "M2"=23.63579766666667; "N2"=70.907393; "O2"=0.06827844587621175; "P2"=0.06827844587621175; "Q2"=52.69300127293044; "R2"=474.2370114563739; "S2"=0.040213437235638; "T2"=0.040213437235638;
    "E3"=3; "G3"=2.229372666666666; "H3"=6.688117999999999; "I3"=0.5889623983027473; "J3"=0.5889623983027473; "K3"=3; "M3"=181.2883913333334; "N3"=543.865174; "O3"=0.5237009467675041; "P3"=0.523700946767504; "Q3"=404.1593844225035; "R3"=3637.434459802532; "S3"=0.3084401656016086; "T3"=0.3084401656016086;
    "E4"=3; "G4"=2.229372666666666; "H4"=6.688117999999999; "I4"=0.5889623983027473; "J4"=0.5889623983027473; "K4"=3; "M4"=111.1005463333333; "N4"=333.301639; "O4"=0.3209442197221123; "P4"=0.3209442197221123; "Q4"=247.6845212472668; "R4"=2229.160691225401; "S4"=0.1890240773689391; "T4"=0.1890240773689391;
    "E5"=3; "G5"=2.229372666666666; "H5"=6.688117999999999; "I5"=0.5889623983027473; "J5"=0.5889623983027473; "K5"=3; "M5"=30.14303933333333; "N5"=90.42911799999999; "O5"=0.08707638763417187; "P5"=0.08707638763417187; "Q5"=67.20006797999154; "R5"=604.8006118199239; "S5"=0.05128471809656156; "T5"=0.05128471809656156;
    "E6"=3; "G6"=0.874264; "H6"=2.622792; "I6"=0.2309657016471988; "J6"=0.2309657016471987; "K6"=3; "M6"=23.63579766666667; "N6"=70.907393; "O6"=0.06827844587621175; "P6"=0.06827844587621175; "Q6"=20.66392701125067; "R6"=185.975343101256; "S6"=0.01576997915917953; "T6"=0.01576997915917953;
    "E7"=3; "G7"=0.874264; "H7"=2.622792; "I7"=0.2309657016471988; "J7"=0.2309657016471987; "K7"=3; "M7"=181.2883913333334; "N7"=543.865174; "O7"=0.5237009467675041; "P7"=0.523700946767504; "Q7"=158.4939141606453; "R7"=1426.445227445808; "S7"=0.1209569566234589; "T7"=0.1209569566234588;
    "E8"=3; "G8"=0.874264; "H8"=2.622792; "I8"=0.2309657016471988; "J8"=0.2309657016471987; "K8"=3; "M8"=111.1005463333333; "N8"=333.301639; "O8"=0.3209442197221123; "P8"=0.3209442197221123; "Q8"=97.13120803956534; "R8"=874.1808723560879; "S8"=0.07412710689773039; "T8"=0.07412710689773037;
    "E9"=3; "G9"=0.874264; "H9"=2.622792; "I9"=0.2309657016471988; "J9"=0.2309657016471987; "K9"=3; "M9"=30.14303933333333; "N9"=90.42911799999999; "O9"=0.08707638763417187; "P9"=0.08707638763417187; "Q9"=26.35297413971733; "R9"=237.176767257456; "S9"=0.02011165896682997; "T9"=0.02011165896682996;
    "E10"=3; "G10"=0.136774; "H10"=0.410322; "I10"=0.0361333680411111; "J10"=0.0361333680411111; "K10"=3; "M10"=23.63579766666667; "N10"=70.907393; "O10"=0.06827844587621175; "P10"=0.06827844587621175; "Q10"=3.232762590060666; "R10"=29.094863310546; "S10"=0.002467130214120244; "T10"=0.002467130214120244;
    "E11"=3; "G11"=0.136774; "H11"=0.410322; "I11"=0.0361333680411111; "J11"=0.0361333680411111; "K11"=3; "M11"=181.2883913333334; "N11"=543.865174; "O11"=0.5237009467675041; "P11"=0.523700946767504; "Q11"=24.79553843622533; "R11"=223.159845926028; "S11"=0.01892307905302856; "T11"=0.01892307905302856;
    "E12"=3; "G12"=0.136774; "H12"=0.410322; "I12"=0.0361333680411111; "J12"=0.0361333680411111; "K12"=3; "M12"=111.1005463333333; "N12"=333.301639; "O12"=0.3209442197221123; "P12"=0.3209442197221123; "Q12"=15.19566612419533; "R12"=136.760995117758; "S12"=0.01159679561188631; "T12"=0.01159679561188631;
    "E13"=3; "G13"=0.136774; "H13"=0.410322; "I13"=0.0361333680411111; "J13"=0.0361333680411111; "K13"=3; "M13"=30.14303933333333; "N13"=90.42911799999999; "O13"=0.08707638763417187; "P13"=0.08707638763417187; "Q13"=4.122784061777332; "R13"=37.10505655599599; "S13"=0.003146363162075988; "T13"=0.003146363162075988;
    "E14"=3; "G14"=0.544844; "H14"=1.634532; "I14"=0.1439385320089428; "J14"=0.1439385320089428; "K14"=3; "M14"=23.63579766666667; "N14"=70.907393; "O14"=0.06827844587621175; "P14"=0.06827844587621175; "Q14"=12.87782254389733; "R14"=115.900402895076; "S14"=0.009827899267273971; "T14"=0.009827899267273971;
    "E15"=3; "G15"=0.544844; "H15"=1.634532; "I15"=0.1439385320089428; "J15"=0.1439385320089428; "K15"=3; "M15"=181.2883913333334; "N15"=543.865174; "O15"=0.5237009467675041; "P15"=0.523700946767504; "Q15"=98.77389228761868; "R15"=888.9650305885681; "S15"=0.07538074548940803; "T15"=0.07538074548940801;
    "E16"=3; "G16"=0.544844; "H16"=1.634532; "I16"=0.1439385320089428; "J16"=0.1439385320089428; "K16"=3; "M16"=111.1005463333333; "N16"=333.301639; "O16"=0.3209442197221123; "P16"=0.3209442197221123; "Q16"=60.53246606643867; "R16"=544.792194597948; "S16"=0.04619623984355642; "T16"=0.04619623984355642;
    "E17"=3; "G17"=0.544844; "H17"=1.634532; "I17"=0.1439385320089428; "J17"=0.1439385320089428; "K17"=3; "M17"=30.14303933333333; "N17"=90.42911799999999; "O17"=0.08707638763417187; "P17"=0.08707638763417187; "Q17"=16.42325412253066; "R17"=147.809287102776; "S17"=0.01253364740870436; "T17"=0.01253364740870436;
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

